$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        $count = $parts.Count
        $revParts = @()
        for ($i = $count - 1; $i -ge 0; $i--) {
            $revParts += $parts[$i]
        }
        $newVal = [string]::Join(", ", $revParts)
        $cell.Value2 = $newVal
    }
}
